# Weekly data refresh: a new price-report row for Jengibre (Vega Central
# Mapocho de Santiago) is inserted at the top of the data block (row 133),
# pushing the previously existing rows 133-142 down to 134-143 (all of their
# data is preserved as-is). The new row carries the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 133; everything from the old row 133
# down shifts one row lower (old 133->134, ..., old 142->143).
$ws.Rows(133).EntireRow.Insert()

# Populate the newly inserted row 133 with this week's record.
$ws.Range("A133").Value = 9
$ws.Range("B133").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C133").Value = "Metropolitana"
$ws.Range("D133").Value = 45166
$ws.Range("E133").Value = 13
$ws.Range("F133").Value = 100114007
$ws.Range("G133").Value = "Jengibre"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 520
$ws.Range("K133").Value = 14000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 14500
$ws.Range("N133").Value = "`$/caja 13 kilos"
$ws.Range("O133").Value = "Perú"
$ws.Range("P133").Value = 1115
$ws.Range("Q133").Value = 13
$ws.Range("R133").Value = "Hortaliza"
